# Apply cryptos list update (commit: "Updated cryptos list on Tue Jul 25 06:46:23 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.136.66"
$ws.Range("E2").Value = "  -2.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.36"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.15"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6867"
$ws.Range("E6").Value = "  -5.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07702"
$ws.Range("E8").Value = "  +7.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3031"
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.12"
$ws.Range("E10").Value = "  -5.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08149"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.854.76"
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7213"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.199"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "88.88"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.159.08"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007798"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.722"
$ws.Range("E18").Value = "  -4.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "233.68"
$ws.Range("E20").Value = "  -5.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.102.01"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.503"
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.80"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.950"
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1419"
$ws.Range("E27").Value = "  -7.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.04"
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.957"
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.400"
$ws.Range("E30").Value = "  -3.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.504"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.484"
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.995"
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05181"
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.173"
$ws.Range("E35").Value = "  -4.40%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7010"
$ws.Range("E36").Value = "  -4.89%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.007"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.650"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01852"
$ws.Range("E39").Value = "  -4.28%  "
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9106"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.093.57"
$ws.Range("E42").Value = "  +5.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.965"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4265"
$ws.Range("E44").Value = "  -4.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.26"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.72"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.751"
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.998.45"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.126"
$ws.Range("E50").Value = "  -4.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.901"
$ws.Range("E51").Value = "  -7.56%  "
